$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1 / A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 12:22"

# Row 18 (Suiza) - updated case counts
$ws.Range("B18").Value = 28496
$ws.Range("C18").Value = 228
$ws.Range("E18").Value = 7087

# Row 35 (Rumania) - updated case counts
$ws.Range("B35").Value = 10096
$ws.Range("C35").Value = 386
$ws.Range("D35").Value = 2478
$ws.Range("E35").Value = 7091
$ws.Range("F35").Value = 236

# Row 56 (Marruecos) - updated case counts
$ws.Range("B56").Value = 3537
$ws.Range("C56").Value = 91
$ws.Range("D56").Value = 430
$ws.Range("E56").Value = 2956
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 151

# Rows 75-76: Bosnia y Herzegovina overtakes Lituania in ranking (swap + update data)
$ws.Range("A75").Value = "Bosnia y Herzegovina"
$ws.Range("B75").Value = 1413
$ws.Range("C75").Value = 45
$ws.Range("D75").Value = 485
$ws.Range("E75").Value = 874
$ws.Range("F75").Value = 4
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 54

$ws.Range("A76").Value = "Lituania"
$ws.Range("B76").Value = 1398
$ws.Range("C76").Value = 28
$ws.Range("D76").Value = 399
$ws.Range("E76").Value = 961
$ws.Range("F76").Value = 17
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 38

# Rows 106-107: Senegal overtakes Malta in ranking (swap + update data)
$ws.Range("A106").Value = "Senegal"
$ws.Range("B106").Value = 479
$ws.Range("C106").Value = 37
$ws.Range("D106").Value = 257
$ws.Range("E106").Value = 216
$ws.Range("F106").Value = 1
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 6

$ws.Range("A107").Value = "Malta"
$ws.Range("B107").Value = 444
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 165
$ws.Range("E107").Value = 276
$ws.Range("F107").Value = 2
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 3
